$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Data edits: C2 becomes text "bro" (was numeric 1), and three new
#    summary rows (Average / Max / Min) are appended below the Total row.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "bro"

$ws.Range("C5").Value = "Average"
$ws.Range("D5").Formula = "=AVERAGE(D2:D2)"

$ws.Range("C6").Value = "Max"
$ws.Range("D6").Formula = "=MAX(D2:D2)"

$ws.Range("C7").Value = "Min"
$ws.Range("D7").Formula = "=MIN(D2:D2)"

# ---------------------------------------------------------------------
# 2. Column widths: split the former uniform B:D width so column C is
#    wider (20.71 chars) and widen the Info column E (40.71 chars).
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.8333333
$ws.Columns.Item(5).ColumnWidth = 39.8333333

# ---------------------------------------------------------------------
# 3. Freeze the header row.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 4. Header row (row 1): blue fill, white bold font, centered.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:E1")
$headerRange.Interior.Color = 12419407
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Date column (B2) centered.
# ---------------------------------------------------------------------
$ws.Range("B2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6. Amount column (D2): drop the "$" from the currency format, keep
#    it as a plain thousands format, right aligned.
# ---------------------------------------------------------------------
$ws.Range("D2").NumberFormat = "#,##0.00"
$ws.Range("D2").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 7. Total row (row 3): orange fill, bold font (keep default color),
#    right aligned.
# ---------------------------------------------------------------------
$totalRange = $ws.Range("C3:D3")
$totalRange.Interior.Color = 49407
$totalRange.Font.Bold = $true
$totalRange.HorizontalAlignment = -4152
$totalRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 8. Average / Max / Min rows (5-7): light-blue fill, bold font, right
#    aligned. Reset the number format the AVERAGE/MAX/MIN formulas
#    auto-inherited from D2 (currency) back to General, and give the
#    brand-new cells the same thin border the rest of the table uses.
# ---------------------------------------------------------------------
$summaryRange = $ws.Range("C5:D7")
$summaryRange.NumberFormat = "General"
$summaryRange.Interior.Color = 15853276
$summaryRange.Font.Bold = $true
$summaryRange.HorizontalAlignment = -4152
$summaryRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 9. AutoFilter strictly over the header + first data row (A1:E2),
#    *not* the Total row. Temporarily clear row 3 so AutoFilter() does
#    not auto-expand into the contiguous block below, then restore it.
# ---------------------------------------------------------------------
$c3Value = $ws.Range("C3").Value2
$d3Formula = $ws.Range("D3").Formula
$ws.Range("C3:D3").ClearContents()
$ws.Range("A1:E2").AutoFilter()
$ws.Range("C3").Value = $c3Value
$ws.Range("D3").Formula = $d3Formula
# Re-entering the formula makes it re-inherit D2's number format; put it
# back to General to match the Total row's original (unformatted) look.
$ws.Range("D3").NumberFormat = "General"

# Excel records the autofilter's backing range as a hidden workbook-level
# defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Income Sheet'!`$A`$1:`$E`$2")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# 10. Conditional formatting: amount greater than 1000 -> green text on
#     green fill (the classic "Good" style).
# ---------------------------------------------------------------------
$cond = $ws.Range("D2").FormatConditions.Add(1, 5, 1000)
$cond.Font.Color = 24832
$cond.Interior.Color = 13561798
